# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$oldText = $cellA1.Value2
$newText = $oldText -replace [regex]::Escape("1000 Bs = 9.87 = 41514.46 pesos"), "1000 Bs = 9.79 = 41143.08 pesos"
$newText = $newText -replace [regex]::Escape("41514.46 pesos = 9.81 = 951.09 Bs"), "41143.08 pesos = 9.78 = 971.25 Bs"
$cellA1.Value2 = $newText

# --- tasas: update the N10/O10/N12/O12 rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value2 = 102.18
$wsTasas.Range("O10").Value2 = 4204
$wsTasas.Range("N12").Value2 = 4209
$wsTasas.Range("O12").Value2 = 99.36
